$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.418.04"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.566.08"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.787.75"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "1.573.78"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "27.414.65"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "1.372.53"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D47").Value = "1.701.04"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -0.86%  "
